$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value. All of these look numeric/percentage
# so we force the cell NumberFormat to Text ("@") before assigning, ensuring Excel
# stores them as literal text strings (matching the source data which is text).
$updates = @{
    'D2' = '307.90'
    'E2' = '1.18%'
    'D3' = '36.67'
    'E3' = '2.65%'
    'D4' = '5.029'
    'E4' = '-0.64%'
    'D5' = '0.07871'
    'E5' = '0.64%'
    'D6' = '2.174'
    'E6' = '-3.72%'
    'D7' = '8.049'
    'E7' = '-0.85%'
    'D8' = '4.065'
    'E8' = '1.61%'
    'D9' = '0.9271'
    'E9' = '0.00%'
    'D10' = '0.09952'
    'E10' = '1.42%'
    'D11' = '0.1885'
    'E11' = '3.62%'
    'D12' = '0.08729'
    'E12' = '0.40%'
    'D13' = '0.03579'
    'E13' = '4.90%'
    'D14' = '0.09947'
    'E14' = '0.23%'
    'D15' = '0.001487'
    'E15' = '0.25%'
    'D16' = '0.005656'
    'E16' = '-0.51%'
    'E17' = '-0.66%'
    'D19' = '0.3434'
    'E19' = '0.06%'
    'D20' = '0.1336'
    'E20' = '1.15%'
    'D21' = '4.943'
    'E21' = '8.68%'
    'D22' = '0.2201'
    'E22' = '-1.54%'
    'D23' = '0.04629'
    'E23' = '-0.86%'
    'D24' = '0.005207'
    'E24' = '15.79%'
    'D25' = '0.001235'
    'E25' = '-0.23%'
    'D26' = '0.0001402'
    'E26' = '7.89%'
    'D27' = '0.0002721'
    'E27' = '0.80%'
    'D39' = '0.01828'
    'E39' = '3.79%'
    'D40' = '0.04771'
    'E40' = '1.45%'
    'D41' = '0.007946'
    'E41' = '-0.55%'
    'D42' = '0.1413'
    'E42' = '-0.43%'
    'D43' = '0.007603'
    'E43' = '-8.18%'
    'D44' = '0.002183'
    'E44' = '-5.05%'
    'E45' = '10.74%'
    'D46' = '0.00006319'
    'E46' = '3.10%'
    'D47' = '0.00000000751'
    'E47' = '0.12%'
    'D49' = '32.38'
    'E49' = '470.57%'
    'D50' = '0.002693'
    'E50' = '0.10%'
    'D51' = '0.00002103'
    'E51' = '0.12%'
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
}
